# Add PF/1.0.5 to meta-sheet
# The sheet holds a small matrix of environment/version labels:
#   Row 1: dev2 | sit2 | uat2 | prod
#   Row 2: PF/1.0.0 | PF/1.0.0 | PF/1.0.0 | PF/1.0.0
# A new release row is appended for PF/1.0.5, marking every environment
# column with "X" since it has not been promoted anywhere yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
